$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$formula = '="---' + "`n" + 'date: " & F4 & "' + "`n" + 'layout: entry' + "`n" + 'name: " & A4 & "' + "`n" + 'description: " & B4 & " ' + "`n" + 'link: " & C4 & "' + "`n" + 'categories:' + "`n" + '" & E4 & "' + "`n" + 'tags:' + "`n" + '" & D4 & "' + "`n" + '---"'

$ws.Range("H4:H9").Formula = $formula

$ws.Range("E16").Select()

$wb.Windows.Item(1).Width = 19200
$wb.Windows.Item(1).Height = 4995
